$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 19.02.2022 16:00"

# Update D3 from text "+0.2" to numeric 0.2
$ws.Range("D3").Value = 0.2

# Update E3 from text date to a real date/time value, matching style of other rows
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("E3").Value = 44611.65723379629
